$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) The intro sentence was originally split into two runs ("W" and
#    "hat is ") with the document's "_GoBack" bookmark sandwiched between
#    them. Re-running a find/replace over that exact text collapses it back
#    into a single "What is " run and drops the now-redundant bookmark
#    (it gets relocated below, in step 2).
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("What is ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "What is ", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2) Append the new "Commands and strategy ..." content after the existing
#    "Pull request:" paragraph, add a blank paragraph as a separator, and
#    mark the final paragraph's paragraph-mark font with the eastAsia hint
#    that Word attached when this text was typed. We build the replacement
#    as a small OOXML fragment and push it in with Range.InsertXML, which
#    is the supported way to splice exact markup (including the relocated
#    "_GoBack" bookmark) into the document.
# ---------------------------------------------------------------------------

function New-Run([string]$text, [bool]$hint) {
    $rPr = '<w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"'
    if ($hint) { $rPr += ' w:hint="eastAsia"' }
    $rPr += '/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr>'
    $preserve = ''
    if ($text -ne $text.Trim() -or $text -eq '') { $preserve = ' xml:space="preserve"' }
    return "<w:r>$rPr<w:t$preserve>$text</w:t></w:r>"
}

# Paragraph mark (pPr) of the "Pull request:" paragraph gains the eastAsia
# hint on its rFonts.
$pullRequestPara = '<w:p w14:paraId="4E0B7751" w14:textId="7EA8E1F8" w:rsidR="00936054" w:rsidRPr="00157186" w:rsidRDefault="00E776C1" w:rsidP="00157186">' + `
    '<w:pPr><w:spacing w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman" w:hint="eastAsia"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr>' + `
    '<w:r w:rsidRPr="00157186"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Pull request:</w:t></w:r>' + `
    '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' + `
    '<w:r w:rsidRPr="00157186"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">Pull requests are proposed changes to a repository submitted by a user and accepted or rejected by a repository''s collaborators. </w:t></w:r>' + `
    '</w:p>'

# New blank separator paragraph (also carries the eastAsia hint, matching
# the formatting Word applied when the cursor moved there).
$blankPara = '<w:p><w:pPr><w:spacing w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman" w:hint="eastAsia"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr></w:p>'

# New "Commands and strategy ..." paragraph, built run-by-run so the
# relocated "_GoBack" bookmark can be dropped exactly between "co" and
# "mmit".
$commandsRuns = `
    (New-Run 'C' $true) + `
    (New-Run 'ommands and strategy ' $false) + `
    (New-Run 'use' $false) + `
    (New-Run 'd' $false) + `
    (New-Run ' ' $false) + `
    (New-Run 'to do' $false) + `
    (New-Run ' part' $true) + `
    (New-Run ' 7: fork' $false) + `
    (New-Run ',' $false) + `
    (New-Run ' co' $false) + `
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' + `
    (New-Run 'mmit' $false) + `
    (New-Run ',' $false) + `
    (New-Run ' pull request' $false) + `
    (New-Run '.' $false)

$commandsPara = '<w:p><w:pPr><w:spacing w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr>' + `
    $commandsRuns + '</w:p>'

$body = $pullRequestPara + $blankPara + $commandsPara

$xmlSnippet = '<?xml version="1.0" encoding="UTF-16" standalone="yes"?>' + `
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
    '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">' + `
    "<w:body>$body</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>"

# Target exactly the final paragraph's contents (excluding its trailing
# paragraph mark) so InsertXML replaces that paragraph in place and appends
# the new paragraphs after it, rather than disturbing anything earlier in
# the document.
$lastPara = $d.Paragraphs.Last
$target = $d.Range($lastPara.Range.Start, $lastPara.Range.End - 1)
$target.InsertXML($xmlSnippet) | Out-Null
